# Rearranges monthly rows so each calendar year's Oct/Nov/Dec figures lead
# that year's block (Jan-Sep follow), matching the source system's fiscal
# re-grouping, and appends newly published 2022 (Oct-Sep rotated) and
# 2023 (Jan-Jul) data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A carries a dedicated style (bordered, centered, top-aligned).
# Grab it from an existing styled cell so newly created rows match.
$styleSrc = $ws.Range("A2")

$data = @(
    ,@("2011-10",99.8034,99.9666,99.9821,100.1865,99.971)
    ,@("2011-11",99.7209,99.8751,99.9288,100.1883,99.8108)
    ,@("2011-12",99.7842,99.9144,99.8379,100.3299,99.8723)
    ,@("2011-01",100.4901,100.5379,99.9218,100.5397,100.9693)
    ,@("2011-02",100.3829,100.3343,99.8724,100.1221,100.6889)
    ,@("2011-03",100.5167,100.399,100.1867,100.5159,100.4091)
    ,@("2011-04",100.302,100.3606,100.0747,100.6168,100.4741)
    ,@("2011-05",100.3152,100.2062,100.129,100.0919,100.232)
    ,@("2011-06",100.1521,100.1423,99.8193,100.1792,100.3254)
    ,@("2011-07",100.192,100.2467,99.9895,100.2018,100.4632)
    ,@("2011-08",100.2502,100.2593,99.9754,100.2827,100.4338)
    ,@("2011-09",100.172,100.2546,99.8217,100.2856,100.5656)
    ,@("2012-10",100.047,100.055,100.1021,100.1269,99.9991)
    ,@("2012-11",99.8177,99.9859,100.0175,100.1423,100.0074)
    ,@("2012-12",99.9629,100.1026,99.9838,100.1402,100.2545)
    ,@("2012-01",100.0121,100.1372,100.209,100.0808,100.1948)
    ,@("2012-02",100.2257,99.9213,99.6433,100.0586,99.8511)
    ,@("2012-03",100.2856,100.1089,99.8693,100.0553,100.1764)
    ,@("2012-04",100.1904,100.0995,99.7108,100.1814,100.2642)
    ,@("2012-05",100.1625,99.9828,99.9636,100.0331,99.8583)
    ,@("2012-06",100.1061,99.9588,99.8553,100.2016,99.8299)
    ,@("2012-07",99.9037,99.9666,99.9915,100.0349,99.962)
    ,@("2012-08",99.9078,100.0122,99.9647,100.0313,100.103)
    ,@("2012-09",100.18,100.1962,99.9146,100.2212,100.3818)
    ,@("2013-10",100.0,100.1,100.0,100.3,100.1)
    ,@("2013-11",100.0,100.1,100.0,100.1,100.2)
    ,@("2013-12",100.0,100.0,99.9,100.0,100.1)
    ,@("2013-01",100.1332,100.2922,99.9976,100.3071,100.5882)
    ,@("2013-02",99.932,99.9798,99.7847,100.0181,100.1253)
    ,@("2013-03",99.927,99.8819,99.9177,99.9339,99.8071)
    ,@("2013-04",100.0107,99.872,99.9073,100.1088,99.6587)
    ,@("2013-05",99.798,99.7836,99.7491,99.9592,99.723)
    ,@("2013-06",99.8549,99.8477,99.8569,100.006,99.7698)
    ,@("2013-07",99.923,99.9731,99.9369,100.0443,99.9998)
    ,@("2013-08",100.0,100.0,99.8,100.1,100.2)
    ,@("2013-09",100.0,100.1,100.0,100.2,100.3)
    ,@("2014-10",99.9,100.0,100.1,100.2,100.0)
    ,@("2014-11",99.8,100.0,100.0,100.1,99.9)
    ,@("2014-12",99.9,100.0,100.0,100.1,99.9)
    ,@("2014-01",100.0,99.9,99.9,99.9,99.9)
    ,@("2014-02",100.1,100.0,100.0,100.1,99.9)
    ,@("2014-03",100.2,100.0,100.0,100.0,99.9)
    ,@("2014-04",100.1,100.0,100.1,100.0,100.0)
    ,@("2014-05",100.0,100.0,99.9,100.0,100.0)
    ,@("2014-06",100.0,99.9,99.6,100.0,100.0)
    ,@("2014-07",100.0,100.0,99.9,100.0,100.0)
    ,@("2014-08",100.0,100.0,99.9,100.0,100.1)
    ,@("2014-09",99.9,100.0,99.9,100.2,100.0)
    ,@("2015-10",100.0,99.9,99.9,100.1,99.7)
    ,@("2015-11",99.869,99.8993,99.7995,100.0082,99.9398)
    ,@("2015-12",99.9243,99.9972,99.9213,100.0636,100.0675)
    ,@("2015-01",99.9,99.9,99.8,100.2,99.9)
    ,@("2015-02",100.0,100.0,100.0,100.0,100.0)
    ,@("2015-03",99.9,100.0,100.1,100.0,100.0)
    ,@("2015-04",99.9,99.9,99.8,99.8,100.1)
    ,@("2015-05",100.0,100.0,99.7,100.0,100.1)
    ,@("2015-06",99.8,100.0,100.0,100.0,100.0)
    ,@("2015-07",99.8,100.0,99.8,100.0,100.1)
    ,@("2015-08",99.9,100.0,99.9,100.0,100.0)
    ,@("2015-09",100.1,100.1,100.1,100.4,99.9)
    ,@("2016-10",99.7,100.0,100.0,100.4,100.0)
    ,@("2016-11",100.4,100.2,100.0,100.2,100.3)
    ,@("2016-12",100.4,100.3,100.1,100.3,100.4)
    ,@("2016-01",100.1462,100.1221,99.8992,100.2538,100.2123)
    ,@("2016-02",100.2002,100.0947,99.9225,99.9982,100.1778)
    ,@("2016-03",100.2138,100.1048,99.859,100.0477,100.221)
    ,@("2016-04",99.9447,99.935,99.7979,99.7926,100.0671)
    ,@("2016-05",100.1,100.0,99.7,99.9,100.0)
    ,@("2016-06",100.1,100.0,99.9,100.2,100.0)
    ,@("2016-07",100.3,100.0,100.1,100.1,99.8)
    ,@("2016-08",100.0,99.9,100.0,99.8,99.9)
    ,@("2016-09",99.7999999999999,100.0,99.9,100.2,100.1)
    ,@("2017-10",100.3,100.1,100.1,100.1,99.9)
    ,@("2017-11",100.2,100.1,99.9,100.1,100.0)
    ,@("2017-12",100.3,100.2,100.0,100.2,100.2)
    ,@("2017-01",100.4,100.2,100.1,100.0,100.3)
    ,@("2017-02",100.2,100.1,99.9,100.2,100.0)
    ,@("2017-03",100.1,100.0,100.1,100.0,99.8)
    ,@("2017-04",99.9,99.9,100.1,100.0,99.8)
    ,@("2017-05",99.8,99.9,100.0,100.0,99.8)
    ,@("2017-06",100.0,99.9,99.8,100.0,99.8)
    ,@("2017-07",99.9,100.0,100.0,100.0,100.1)
    ,@("2017-08",100.2,100.1,99.9,100.0,100.3)
    ,@("2017-09",100.2,100.1,99.9,100.1,100.2)
    ,@("2018-10",100.2,100.0,99.8,100.3,99.9)
    ,@("2018-11",100.1,100.2,100.1,100.4,100.2)
    ,@("2018-12",99.9,100.0,100.0,100.2,100.1)
    ,@("2018-01",100.1,100.1,100.0,100.2,100.0)
    ,@("2018-02",99.9,100.0,100.1,100.0,100.0)
    ,@("2018-03",99.9,99.8,99.9,99.8,99.8)
    ,@("2018-04",100.0,99.9,99.8,99.8,99.9)
    ,@("2018-05",100.0,100.0,99.9,100.1,100.1)
    ,@("2018-06",100.0,100.0,100.0,99.9,100.1)
    ,@("2018-07",99.9,100.2,100.2,100.4,100.2)
    ,@("2018-08",100.3,100.3,100.2,100.4,100.3)
    ,@("2018-09",100.1,100.1,99.9,100.0,100.3)
    ,@("2019-10",99.8,100.3,99.9,100.0,101.1)
    ,@("2019-11",99.7,100.3,99.8,100.0,101.1)
    ,@("2019-12",99.9,99.8,99.7,99.8,99.7)
    ,@("2019-01",100.0,100.0,99.9,100.2,99.9)
    ,@("2019-02",99.9,99.9,99.5,99.9,100.1)
    ,@("2019-03",99.9,100.0,99.8,99.9,100.2)
    ,@("2019-04",100.0,100.2,99.8,99.9,100.6)
    ,@("2019-05",100.1,100.1,99.8,100.0,100.4)
    ,@("2019-06",100.1,100.0,100.0,100.0,100.0)
    ,@("2019-07",100.2,100.0,99.9,100.1,100.0)
    ,@("2019-08",100.1,100.2,99.5,100.0,100.8)
    ,@("2019-09",100.3,100.5,100.1,100.1,101.0)
    ,@("2020-10",99.9,99.9,100.0,99.9,99.8)
    ,@("2020-11",100.0,100.1,99.8,100.0,100.3)
    ,@("2020-12",100.0,100.1,99.7,99.8,100.5)
    ,@("2020-01",99.9,100.0,99.9,100.0,100.0)
    ,@("2020-02",100.1,100.0,99.8,99.9,100.0)
    ,@("2020-03",99.8,99.8,99.9,99.7,99.7)
    ,@("2020-04",99.9,99.9,99.9,99.8,99.9)
    ,@("2020-05",99.9,99.7,99.8,99.9,99.5)
    ,@("2020-06",100.0,100.1,99.9,99.9,100.3)
    ,@("2020-07",99.8,100.1,100.0,99.7,100.6)
    ,@("2020-08",100.5,100.1,99.6,99.8,100.2)
    ,@("2020-09",99.9,99.9,99.7,99.8,100.0)
    ,@("2021-10",100.3,100.1,99.8,100.3,100.2)
    ,@("2021-11",100.4,100.4,99.8,100.4,100.8)
    ,@("2021-12",100.1,100.0,99.9,99.7,100.1)
    ,@("2021-01",100.1,100.2,99.9,99.8,100.7)
    ,@("2021-02",100.0,100.0,99.8,100.0,100.2)
    ,@("2021-03",100.0,100.2,100.3,100.3,100.2)
    ,@("2021-04",100.3,100.1,100.4,100.2,99.7)
    ,@("2021-05",100.4,100.1,100.1,100.0,100.0)
    ,@("2021-06",99.8,99.8,100.1,99.7,99.5)
    ,@("2021-07",100.0,100.1,100.2,100.2,100.0)
    ,@("2021-08",100.0,100.0,99.8,100.3,100.0)
    ,@("2021-09",99.9,100.0,100.0,100.2,99.9)
    ,@("2022-10",100.3,100.5,100.7,100.5,100.6)
    ,@("2022-11",100.2,100.1,100.0,100.2,100.2)
    ,@("2022-12",100.1,99.8,100.1,99.8,99.4)
    ,@("2022-01",100.1,100.0,100.2,100.2,99.9)
    ,@("2022-02",100.2,100.1,99.9,99.9,100.2)
    ,@("2022-03",100.4,100.2,100.0,99.8,100.5)
    ,@("2022-04",100.2,100.2,100.0,100.0,100.4)
    ,@("2022-05",100.2,100.3,100.1,100.5,100.6)
    ,@("2022-06",100.1,100.3,99.9,100.7,100.5)
    ,@("2022-07",99.7,100.2,100.2,100.1,100.6)
    ,@("2022-08",100.0,99.9,99.6,100.3,100.0)
    ,@("2022-09",100.0,100.1,100.0,100.3,100.3)
    ,@("2023-01",100.0,99.7,99.8,99.7,99.5)
    ,@("2023-02",99.7,99.7,99.5,99.7,99.9)
    ,@("2023-03",100.2,100.0,99.9,100.1,99.9)
    ,@("2023-04",100.0,99.7,99.7,100.0,99.5)
    ,@("2023-05",100.1,99.8,99.6,99.9,99.7)
    ,@("2023-06",100.1,99.8,99.6,100.3,99.8)
    ,@("2023-07",100.2,100.3,100.2,100.6,100.3)
)

$startRow = 2
$lastExistingRow = 133

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    if ($r -gt $lastExistingRow) {
        # Brand-new row past the old A1:F133 range - clone column A's style
        # (border + centered/top alignment) before writing its value.
        $dst = $ws.Range("A$r")
        $styleSrc.Copy($dst)
    }

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

Write-Host "Rewrote $($data.Count) data rows (rows $startRow..$($startRow + $data.Count - 1))"
